$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22, shifting existing rows 22:125 down to 23:126
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record
$ws.Cells.Item(22, 1).Value = 3
$ws.Cells.Item(22, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44561
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = 100112026
$ws.Cells.Item(22, 7).Value = "Haba"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 40
$ws.Cells.Item(22, 11).Value = 8000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 8000
$ws.Cells.Item(22, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(22, 16).Value = 320
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
